$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new data row (row 5): ID, NAME, EMAIL, FILE
$ws.Range("A5").Value = "123"
$ws.Range("A5").NumberFormat = $ws.Range("A4").NumberFormat

$ws.Range("B5").Value = "Glaiza Marie"

$ws.Range("C5").Value = "glaiza.garay@powersource.group"
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:glaiza.garay@powersource.group")
$ws.Range("C5").Style = "Hyperlink"

$ws.Range("D5").Value = "file 1.pdf"

# Match the new active cell / selection
$ws.Range("D5").Select()
